{"js": "// Replace the date and the three-digit-by-one-digit multiplication\n// problems throughout the document body (including inside the table).\nconst replacements = [\n  [\"2025-01-24 Friday\", \"2025-01-25 Saturday\"],\n  [\"996\\u00D73=\", \"722\\u00D76=\"],\n  [\"852\\u00D74=\", \"546\\u00D76=\"],\n  [\"357\\u00D72=\", \"720\\u00D79=\"],\n  [\"848\\u00D78=\", \"912\\u00D75=\"],\n  [\"280\\u00D79=\", \"617\\u00D75=\"],\n  [\"954\\u00D75=\", \"742\\u00D77=\"],\n  [\"312\\u00D78=\", \"897\\u00D76=\"],\n  [\"943\\u00D72=\", \"219\\u00D73=\"],\n  [\"401\\u00D75=\", \"876\\u00D77=\"],\n  [\"967\\u00D78=\", \"140\\u00D75=\"],\n  [\"992\\u00D76=\", \"814\\u00D73=\"],\n  [\"855\\u00D72=\", \"868\\u00D74=\"],\n  [\"599\\u00D74=\", \"983\\u00D76=\"],\n  [\"924\\u00D74=\", \"562\\u00D75=\"],\n  [\"451\\u00D78=\", \"363\\u00D72=\"],\n  [\"558\\u00D75=\", \"524\\u00D75=\"],\n  [\"712\\u00D74=\", \"211\\u00D79=\"],\n  [\"258\\u00D73=\", \"492\\u00D74=\"],\n  [\"305\\u00D78=\", \"247\\u00D75=\"],\n  [\"479\\u00D79=\", \"460\\u00D74=\"],\n  [\"419\\u00D76=\", \"556\\u00D72=\"],\n  [\"871\\u00D72=\", \"429\\u00D75=\"],\n  [\"208\\u00D78=\", \"593\\u00D72=\"],\n  [\"886\\u00D75=\", \"929\\u00D77=\"],\n  [\"755\\u00D77=\", \"457\\u00D74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and the three-digit-by-one-digit multiplication\n# problems throughout the document body (including inside the table).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-01-24 Friday\", \"2025-01-25 Saturday\"),\n    @(\"996\u00d73=\", \"722\u00d76=\"),\n    @(\"852\u00d74=\", \"546\u00d76=\"),\n    @(\"357\u00d72=\", \"720\u00d79=\"),\n    @(\"848\u00d78=\", \"912\u00d75=\"),\n    @(\"280\u00d79=\", \"617\u00d75=\"),\n    @(\"954\u00d75=\", \"742\u00d77=\"),\n    @(\"312\u00d78=\", \"897\u00d76=\"),\n    @(\"943\u00d72=\", \"219\u00d73=\"),\n    @(\"401\u00d75=\", \"876\u00d77=\"),\n    @(\"967\u00d78=\", \"140\u00d75=\"),\n    @(\"992\u00d76=\", \"814\u00d73=\"),\n    @(\"855\u00d72=\", \"868\u00d74=\"),\n    @(\"599\u00d74=\", \"983\u00d76=\"),\n    @(\"924\u00d74=\", \"562\u00d75=\"),\n    @(\"451\u00d78=\", \"363\u00d72=\"),\n    @(\"558\u00d75=\", \"524\u00d75=\"),\n    @(\"712\u00d74=\", \"211\u00d79=\"),\n    @(\"258\u00d73=\", \"492\u00d74=\"),\n    @(\"305\u00d78=\", \"247\u00d75=\"),\n    @(\"479\u00d79=\", \"460\u00d74=\"),\n    @(\"419\u00d76=\", \"556\u00d72=\"),\n    @(\"871\u00d72=\", \"429\u00d75=\"),\n    @(\"208\u00d78=\", \"593\u00d72=\"),\n    @(\"886\u00d75=\", \"929\u00d77=\"),\n    @(\"755\u00d77=\", \"457\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
